# Updated cryptos list on Sun Oct 27 08:55:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) and Volume(1h) (column E) updates ---
$ws.Range("D2").Value = "67.193.97"
$ws.Range("D3").Value = "2.476.74"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D5").Value = "584.76"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "174.72"
$ws.Range("E6").Value = "  +3.76%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +2.61%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D14").Value = "25.50"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "67.099.23"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "2.522.84"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "7.57"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "350.15"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "69.11"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "4.24"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "503.09"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "7.76"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "1.76"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "161.15"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "18.15"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("D45").Value = "142.32"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("E46").Value = "  +0.76%  "

# --- Rows 47 and 48 swapped places (ARBITRUM now ranked above BabyDogeCoin) ---
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.514"
$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0" + [char]0x2086 + "0257"
$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").Value = "0.0740"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  +0.57%  "
